# adding profits to tables
# Insert a new "M_PL" (profits) column right after the M_ETR column (B),
# shifting the existing GFA/IMF/OECD columns one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column before column C - existing C:J shift to D:K
$ws.Range("C1").EntireColumn.Insert()

# Header for the new column
$ws.Range("C1").Value = "M_PL"

# Profit values for each group row
$ws.Range("C2").Value = 1007534436142
$ws.Range("C3").Value = -269766813
$ws.Range("C4").Value = 20228669958
$ws.Range("C5").Value = 344467447608
$ws.Range("C6").Value = 45733381438
